$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Duplicate the "LED" pin-assignment block (rows 38:47) down to rows 49:58
#    so the sheet records both the old and the new GPIO_1 pin settings.
# ---------------------------------------------------------------------------
$src = $ws.Range("A38:H47")
$dst = $ws.Range("A49")
$src.Copy($dst)

# Remove the G/H cells that should not exist on the "plain" rows (50:56) -
# only rows 49 (new-settings note) and 57:58 (which mirror 46:47's extra
# styled-but-empty G/H cells) keep them. Row 49 only keeps G (no H).
$ws.Range("G50:H56").Clear()
$ws.Range("H49").Clear()

# ---------------------------------------------------------------------------
# 2. Re-point the "variable" (column D) assignments on the new rows: the new
#    GPIO_1 settings use the previous block's values in reverse order.
# ---------------------------------------------------------------------------
$ws.Cells.Item(49, 4).Value = "LEDAux2"
$ws.Cells.Item(50, 4).Value = "LEDAux1"
$ws.Cells.Item(51, 4).Value = "LED854_Freq"
$ws.Cells.Item(52, 4).Value = "LED854_Power"
$ws.Cells.Item(53, 4).Value = "LED729_RF2"
$ws.Cells.Item(54, 4).Value = "LED729_RF1"
$ws.Cells.Item(55, 4).Value = "LED854"
$ws.Cells.Item(56, 4).Value = "LED729"
$ws.Cells.Item(57, 4).Value = "LED397_2"
$ws.Cells.Item(58, 4).Value = "LED397_1"

# ---------------------------------------------------------------------------
# 3. Annotate which block is the new configuration and which is obsolete.
#    (Insert the "new settings" note first so it lands before the "old
#    settings" note in the shared-string table, matching authoring order.)
# ---------------------------------------------------------------------------
$ws.Cells.Item(49, 7).Value = "These are the new settings"
$ws.Cells.Item(38, 7).Value = "These are the old, now obsolete settings"

# ---------------------------------------------------------------------------
# 4. Unhide columns F:I (pin/old-new comparison helper columns) which used
#    to be hidden, and restore column I to a normal (non-zero) width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).Hidden = $false
$ws.Columns.Item(7).Hidden = $false
$ws.Columns.Item(8).Hidden = $false
$ws.Columns.Item(9).Hidden = $false
$ws.Columns.Item(9).ColumnWidth = 8.25

# ---------------------------------------------------------------------------
# 5. Update the view: scroll the frozen pane down and move the selection to
#    the newly annotated cell.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G38").Select()
